$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Run Number" (column A) values for the two trailing groups of
# rows (62-65 => run 13, 67-70 => run 14) that were left blank.
foreach ($r in 62..65) {
    $ws.Cells.Item($r, 1).Value = 13
}

foreach ($r in 67..70) {
    $ws.Cells.Item($r, 1).Value = 14
}

# Update the active selection to match the recorded view state.
$ws.Range("I64").Select()
